$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original row 13 ("6270264 - Juan Fernando Zapata Zapata" in B13/C13, with no
# A13 label) is removed entirely. Deleting the whole row shifts every following row
# up by one, which already reproduces the vast majority of the target layout
# (labels in column A, and the row heights) without any further work.
$ws.Rows.Item(13).Delete()

# After the shift, a handful of B/C (data) cells end up holding the wrong text and
# need to be corrected to match the new content.

# Row 10 ("Objetivos:") now shows the docente's identification instead of the
# (now removed) objectives paragraph.
$ws.Cells.Item(10, 2).Value = "6270264 - Juan Fernando Zapata Zapata"
$ws.Cells.Item(10, 3).Value = "6270264 - Juan Fernando Zapata Zapata"

# Row 13 ("Programa resumido:") now just says "Semestral".
$ws.Cells.Item(13, 2).Value = "Semestral"
$ws.Cells.Item(13, 3).Value = "Semestral"

# Row 15 ("Programa:") now shows the activation date. Copy it from the existing
# "Ativação:" row (row 8) instead of assigning the literal string, so Excel keeps
# storing it as text (sharing the existing string) instead of auto-converting the
# dd/mm/yyyy-looking text into a date serial number with a brand new number format.
$ws.Cells.Item(8, 2).Copy($ws.Cells.Item(15, 2))
$ws.Cells.Item(8, 3).Copy($ws.Cells.Item(15, 3))

# Row 18 ("Método:") now shows the docente's identification again.
$ws.Cells.Item(18, 2).Value = "6270264 - Juan Fernando Zapata Zapata"
$ws.Cells.Item(18, 3).Value = "6270264 - Juan Fernando Zapata Zapata"

# Row 19 ("Critério:") now shows the evaluation-method paragraph.
$ws.Cells.Item(19, 2).Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$ws.Cells.Item(19, 3).Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."

# Row 20 ("Norma de recuperação:") now shows the passing-grade criterion.
$ws.Cells.Item(20, 2).Value = "NF≥ 5,0."
$ws.Cells.Item(20, 3).Value = "NF≥ 5,0."

# Row 21 ("Bibliografia:") now shows the recovery-grade norm.
$ws.Cells.Item(21, 2).Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
$ws.Cells.Item(21, 3).Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
